# Auto-generated: apply row permutation per diff (rows 2-9 of sheet "Artfynd")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111739317
$ws.Range("B2").Value = 78579
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("Q2").Value = 573911.5177193542
$ws.Range("R2").Value = 7172648.020174325
$ws.Range("A3").Value = 111739311
$ws.Range("B3").Value = 77515
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 574011.8892867711
$ws.Range("R3").Value = 7172473.089384713
$ws.Range("A4").Value = 111739306
$ws.Range("B4").Value = 56398
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 573906.0397215446
$ws.Range("R4").Value = 7172521.061635921
$ws.Range("A5").Value = 111739307
$ws.Range("B5").Value = 56543
$ws.Range("E5").Value = 103021
$ws.Range("F5").Value = "Talltita"
$ws.Range("G5").Value = "Poecile montanus"
$ws.Range("H5").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "3"
$ws.Range("M5").Value = "födosökande"
$ws.Range("Q5").Value = 573960.5743707293
$ws.Range("R5").Value = 7172501.399265604
$ws.Range("A6").Value = 111739315
$ws.Range("B6").Value = 78605
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6462
$ws.Range("F6").Value = "Stuplav"
$ws.Range("G6").Value = "Nephroma bellum"
$ws.Range("H6").Value = "(Spreng.) Tuck."
$ws.Range("A7").Value = 111739313
$ws.Range("B7").Value = 73701
$ws.Range("E7").Value = 1467
$ws.Range("F7").Value = "Rödbrun blekspik"
$ws.Range("G7").Value = "Sclerophora coniophaea"
$ws.Range("H7").Value = "(Norman) J.Mattsson & Middelb."
$ws.Range("Q7").Value = 574025.0565134182
$ws.Range("R7").Value = 7172443.417908707
$ws.Range("A8").Value = 111739316
$ws.Range("B8").Value = 78578
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = "Lunglav"
$ws.Range("G8").Value = "Lobaria pulmonaria"
$ws.Range("H8").Value = "(L.) Hoffm."
$ws.Range("Q8").Value = 573904.5013778479
$ws.Range("R8").Value = 7172636.708955797
$ws.Range("A9").Value = 111739309
$ws.Range("B9").Value = 78536
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 229497
$ws.Range("F9").Value = "Korallblylav"
$ws.Range("G9").Value = "Parmeliella triptophylla"
$ws.Range("H9").Value = "(Ach.) Müll.Arg."
$ws.Range("I9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("Q9").Value = 574011.1276117128
$ws.Range("R9").Value = 7172434.078971106
